$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: add new value in E23
$ws.Range("E23").Value = "falta encontrar ese campo"

# Row 27: add new value in D27
$ws.Range("D27").Value = "CUSTITEMNUMBER_UNI_CANT_BULTO"

# Row 28: add new value in D28
$ws.Range("D28").Value = "CUSTITEMNUMBER_UNI_FRAC_BULTO"

# Row 29: add new value in D29
$ws.Range("D29").Value = "CUSTITEMNUMBER_UNI_NRO_BULTO"

# Row 25: change D25 from "expirationDate" to "expirationdate"
$ws.Range("D25").Value = "expirationdate"

# Update the selected cell to match the final cursor position
$ws.Range("F26").Select()
